# SCD0011 until SCD0016 update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from SCD0184 to SCD0011
$ws.Name = "SCD0011"

# Update the TC_ID cell (B2) from "DGS-199" to "SCD0011-015"
$ws.Range("B2").Value = "SCD0011-015"

# Column B needs to widen to fit the new, longer TC_ID text
$ws.Columns.Item(2).ColumnWidth = 11.6

# Move the active selection to B3 (top-left cell returns to A1, no more D1 freeze/scroll)
$ws.Range("B3").Select() | Out-Null
